$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 102
$ws.Range("D2").Value = 102
$ws.Range("E2").Value = 86
$ws.Range("F2").Value = 0.8431372549019608
$ws.Range("G2").Value = 0.8431372549019608
$ws.Range("H2").Value = 0.09853558668981278
$ws.Range("I2").Value = 0.08307902407180294
$ws.Range("J2").Value = 453398.1905958019
$ws.Range("K2").Value = 164085.0954989009
$ws.Range("M2").Value = 164085.0954989009
$ws.Range("N2").Value = 617483.2860947028
$ws.Range("O2").Value = 10034971.8888
$ws.Range("P2").Value = 9627230.958700001
$ws.Range("Q2").Value = 0.01635132587486724
$ws.Range("R2").Value = 0.01704385157090465

# Row 3
$ws.Range("C3").Value = 103
$ws.Range("D3").Value = 103
$ws.Range("E3").Value = 85
$ws.Range("F3").Value = 0.8252427184466019
$ws.Range("G3").Value = 0.8252427184466019
$ws.Range("H3").Value = 0.09946524789627358
$ws.Range("I3").Value = 0.08208297156488596
$ws.Range("J3").Value = 475130.6665414795
$ws.Range("K3").Value = 172352.9134777698
$ws.Range("M3").Value = 172352.9134777698
$ws.Range("N3").Value = 647483.5800192493
$ws.Range("O3").Value = 10633646.086764
$ws.Range("P3").Value = 10226272.928761
$ws.Range("Q3").Value = 0.01620826121835127
$ws.Range("R3").Value = 0.01685393248140618

# Row 4
$ws.Range("C4").Value = 104
$ws.Range("D4").Value = 104
$ws.Range("E4").Value = 87
$ws.Range("F4").Value = 0.8365384615384616
$ws.Range("G4").Value = 0.8365384615384616
$ws.Range("H4").Value = 0.09774377642832541
$ws.Range("I4").Value = 0.08176642835831067
$ws.Range("J4").Value = 502966.1150419703
$ws.Range("K4").Value = 179184.265134226
$ws.Range("M4").Value = 179184.265134226
$ws.Range("N4").Value = 682150.3801761963
$ws.Range("O4").Value = 10970666.81506692
$ws.Range("P4").Value = 10562322.46232383
$ws.Range("Q4").Value = 0.01633303318337382
$ws.Range("R4").Value = 0.01696447592595118

# Row 5
$ws.Range("C5").Value = 105
$ws.Range("D5").Value = 104
$ws.Range("E5").Value = 87
$ws.Range("F5").Value = 0.8365384615384616
$ws.Range("G5").Value = 0.8285714285714286
$ws.Range("H5").Value = 0.09774030428156209
$ws.Range("I5").Value = 0.08098482354758002
$ws.Range("J5").Value = 516886.2289629109
$ws.Range("K5").Value = 183570.3583230936
$ws.Range("M5").Value = 183570.3583230936
$ws.Range("N5").Value = 700456.5872860046
$ws.Range("O5").Value = 11269370.82291893
$ws.Range("P5").Value = 10858676.13959355
$ws.Range("Q5").Value = 0.01628931740801003
$ws.Range("R5").Value = 0.01690540872231639

# Row 6
$ws.Range("C6").Value = 106
$ws.Range("D6").Value = 106
$ws.Range("E6").Value = 90
$ws.Range("F6").Value = 0.8490566037735849
$ws.Range("G6").Value = 0.8490566037735849
$ws.Range("H6").Value = 0.09459020222215261
$ws.Range("I6").Value = 0.08031243584899751
$ws.Range("J6").Value = 533100.1538977289
$ws.Range("K6").Value = 188616.1381057517
$ws.Range("M6").Value = 188616.1381057517
$ws.Range("N6").Value = 721716.2920034805
$ws.Range("O6").Value = 11647629.9738065
$ws.Range("P6").Value = 11233164.44998135
$ws.Range("Q6").Value = 0.01619352078748352
$ws.Range("R6").Value = 0.0167910065721565

# Row 7
$ws.Range("C7").Value = 102
$ws.Range("D7").Value = 101
$ws.Range("E7").Value = 85
$ws.Range("F7").Value = 0.8415841584158416
$ws.Range("G7").Value = 0.8333333333333334
$ws.Range("H7").Value = 0.1003363526021026
$ws.Range("I7").Value = 0.08361362716841886
$ws.Range("J7").Value = 460561.0260389551
$ws.Range("K7").Value = 167666.5132204776
$ws.Range("M7").Value = 167666.5132204776
$ws.Range("N7").Value = 628227.5392594326
$ws.Range("O7").Value = 10081246.3288
$ws.Range("P7").Value = 9673505.398699997
$ws.Range("Q7").Value = 0.01663152627681457
$ws.Range("R7").Value = 0.01733254971284866

# Row 8
$ws.Range("C8").Value = 103
$ws.Range("D8").Value = 103
$ws.Range("E8").Value = 87
$ws.Range("F8").Value = 0.8446601941747572
$ws.Range("G8").Value = 0.8446601941747572
$ws.Range("H8").Value = 0.09903479425028895
$ws.Range("I8").Value = 0.08365074854150621
$ws.Range("J8").Value = 486903.624841487
$ws.Range("K8").Value = 178239.3926277735
$ws.Range("M8").Value = 178239.3926277735
$ws.Range("N8").Value = 665143.0174692603
$ws.Range("O8").Value = 10637203.005464
$ws.Range("P8").Value = 10229829.847461
$ws.Range("Q8").Value = 0.01675622741581762
$ws.Range("R8").Value = 0.01742349533526324

# Row 9
$ws.Range("C9").Value = 104
$ws.Range("D9").Value = 103
$ws.Range("E9").Value = 88
$ws.Range("F9").Value = 0.8543689320388349
$ws.Range("G9").Value = 0.8461538461538461
$ws.Range("H9").Value = 0.09855902998287509
$ws.Range("I9").Value = 0.08339610229320202
$ws.Range("J9").Value = 508874.7839486722
$ws.Range("K9").Value = 182138.599587577
$ws.Range("M9").Value = 182138.599587577
$ws.Range("N9").Value = 691013.3835362492
$ws.Range("O9").Value = 10795112.11552792
$ws.Range("P9").Value = 10386767.76278483
$ws.Range("Q9").Value = 0.01687232125413361
$ws.Range("R9").Value = 0.01753563801052419

# Row 10
$ws.Range("C10").Value = 105
$ws.Range("D10").Value = 105
$ws.Range("E10").Value = 89
$ws.Range("F10").Value = 0.8476190476190476
$ws.Range("G10").Value = 0.8476190476190476
$ws.Range("H10").Value = 0.09827781108665552
$ws.Range("I10").Value = 0.08330214463535565
$ws.Range("J10").Value = 531695.2573335718
$ws.Range("K10").Value = 190974.872508424
$ws.Range("M10").Value = 190974.872508424
$ws.Range("N10").Value = 722670.1298419957
$ws.Range("O10").Value = 11349892.00119376
$ws.Range("P10").Value = 10939197.31786837
$ws.Range("Q10").Value = 0.01682614006268409
$ws.Range("R10").Value = 0.01745785060449367

# Row 11
$ws.Range("C11").Value = 106
$ws.Range("D11").Value = 106
$ws.Range("E11").Value = 89
$ws.Range("F11").Value = 0.839622641509434
$ws.Range("G11").Value = 0.839622641509434
$ws.Range("H11").Value = 0.0979325194061124
$ws.Range("I11").Value = 0.082226360633434
$ws.Range("J11").Value = 547947.509474281
$ws.Range("K11").Value = 196039.8158940278
$ws.Range("M11").Value = 196039.8158940278
$ws.Range("N11").Value = 743987.3253683088
$ws.Range("O11").Value = 11668294.54402957
$ws.Range("P11").Value = 11253829.02020442
$ws.Range("Q11").Value = 0.016801068498424
$ws.Range("R11").Value = 0.01741983244476792

# Row 12
$ws.Range("C12").Value = 102
$ws.Range("D12").Value = 102
$ws.Range("E12").Value = 86
$ws.Range("F12").Value = 0.8431372549019608
$ws.Range("G12").Value = 0.8431372549019608
$ws.Range("H12").Value = 0.09930466255509107
$ws.Range("I12").Value = 0.08372746058566501
$ws.Range("J12").Value = 458202.2395800996
$ws.Range("K12").Value = 166487.1199910498
$ws.Range("M12").Value = 166487.1199910498
$ws.Range("N12").Value = 624689.3595711493
$ws.Range("O12").Value = 10064889.9988
$ws.Range("P12").Value = 9657149.068700001
$ws.Range("Q12").Value = 0.01654137501859429
$ws.Range("R12").Value = 0.01723977944284353

# Row 13
$ws.Range("C13").Value = 103
$ws.Range("D13").Value = 103
$ws.Range("E13").Value = 87
$ws.Range("F13").Value = 0.8446601941747572
$ws.Range("G13").Value = 0.8446601941747572
$ws.Range("H13").Value = 0.1058451956245466
$ws.Range("I13").Value = 0.08940322348869467
$ws.Range("J13").Value = 544119.8952699812
$ws.Range("K13").Value = 206847.5278420206
$ws.Range("M13").Value = 206847.5278420206
$ws.Range("N13").Value = 750967.4231120017
$ws.Range("O13").Value = 10559579.354464
$ws.Range("P13").Value = 10152206.196461
$ws.Range("Q13").Value = 0.01958861436602368
$ws.Range("R13").Value = 0.02037463816624671

# Row 14
$ws.Range("C14").Value = 104
$ws.Range("D14").Value = 104
$ws.Range("E14").Value = 86
$ws.Range("F14").Value = 0.8269230769230769
$ws.Range("G14").Value = 0.8269230769230769
$ws.Range("H14").Value = 0.1118823209942627
$ws.Range("I14").Value = 0.09251807312987113
$ws.Range("J14").Value = 618419.5048371302
$ws.Range("K14").Value = 236910.960031806
$ws.Range("M14").Value = 236910.960031806
$ws.Range("N14").Value = 855330.4648689362
$ws.Range("O14").Value = 10924825.81249792
$ws.Range("P14").Value = 10516481.45975483
$ws.Range("Q14").Value = 0.02168555948606353
$ws.Range("R14").Value = 0.02252758785706347

# Row 15
$ws.Range("C15").Value = 105
$ws.Range("D15").Value = 104
$ws.Range("E15").Value = 87
$ws.Range("F15").Value = 0.8365384615384616
$ws.Range("G15").Value = 0.8285714285714286
$ws.Range("H15").Value = 0.1152466813322768
$ws.Range("I15").Value = 0.09549010738960079
$ws.Range("J15").Value = 674928.3527059811
$ws.Range("K15").Value = 262591.4201946286
$ws.Range("M15").Value = 262591.4201946286
$ws.Range("N15").Value = 937519.7729006096
$ws.Range("O15").Value = 11163372.94887286
$ws.Range("P15").Value = 10752678.26554748
$ws.Range("Q15").Value = 0.0235225877875147
$ws.Range("R15").Value = 0.02442102457729016

# Row 16
$ws.Range("C16").Value = 106
$ws.Range("D16").Value = 106
$ws.Range("E16").Value = 89
$ws.Range("F16").Value = 0.839622641509434
$ws.Range("G16").Value = 0.839622641509434
$ws.Range("H16").Value = 0.1136698443269339
$ws.Range("I16").Value = 0.09543977495374641
$ws.Range("J16").Value = 706045.9233144443
$ws.Range("K16").Value = 275089.0228141093
$ws.Range("M16").Value = 275089.0228141093
$ws.Range("N16").Value = 981134.9461285535
$ws.Range("O16").Value = 11682335.88493904
$ws.Range("P16").Value = 11267870.3611139
$ws.Range("Q16").Value = 0.02354743310956811
$ws.Range("R16").Value = 0.02441357718877013

